$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9 - this shifts the old rows 9-20 down to 10-21,
# keeping row 8 where it is.
$ws.Rows(9).Insert()

# --- Row 8 (unchanged position, only the "Qty executed upto date" count changes) ---
$ws.Cells.Item(8, 3).Value = 69

# --- Row 9 (brand-new row: "Short point (up to 3 mtr.)") ---
$ws.Cells.Item(9, 1).Value = "P. point"
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 82
$ws.Cells.Item(9, 4).Value = "'2"
$ws.Cells.Item(9, 5).Value = "Short point (up to 3 mtr.)"
$ws.Cells.Item(9, 6).Value = 256
$ws.Cells.Item(9, 7).Value = "'20992.00"
$ws.Cells.Item(9, 8).Value = 0

# --- Row 10 (was row 9: Medium point) ---
$ws.Cells.Item(10, 3).Value = 81
$ws.Cells.Item(10, 7).Value = "'38232.00"

# --- Row 11 (was row 10: Long point) ---
$ws.Cells.Item(11, 3).Value = 33
$ws.Cells.Item(11, 7).Value = "'21846.00"

# --- Row 12 (was row 11: Rewiring 3/5 pin plug point) ---
$ws.Cells.Item(12, 3).Value = 98

# --- Row 13 (was row 12: On board) ---
$ws.Cells.Item(13, 3).Value = 56
$ws.Cells.Item(13, 7).Value = "'7616.00"

# --- Row 14 (was row 13: P & F switch) ---
$ws.Cells.Item(14, 3).Value = 89
$ws.Cells.Item(14, 7).Value = "'2047.00"

# --- Row 15 (was row 14: Total) ---
$ws.Cells.Item(15, 3).Value = 85

# --- Row 16 (was row 15: Add Tender Premium) ---
$ws.Cells.Item(16, 3).Value = 59

# --- Row 17 (was row 16: Grand Total) ---
$ws.Cells.Item(17, 3).Value = 53

# --- Row 19 (was row 18: Grand Total Rs.) ---
$ws.Cells.Item(19, 7).Value = "'90733.00"
$ws.Cells.Item(19, 8).Value = "'90733.00"

# --- Row 21 (was row 20: NET PAYABLE AMOUNT Rs.) ---
$ws.Cells.Item(21, 7).Value = "'90733.00"
$ws.Cells.Item(21, 8).Value = "'90733.00"
